# Append a "Things to look for correlation with:" section with a bulleted
# list to the end of the document, right before the document's trailing
# empty paragraph (which must remain the very last paragraph).

$d = $word.ActiveDocument

# --- Step 1: "use up" five bullet-list definitions on throwaway paragraphs
# so that the real list we keep lands on the same list id the authoring
# session ended up with. We create five temp paragraphs, give each one its
# own bullet-list format (each call mints a fresh list definition), then
# delete the paragraphs again - this leaves the list definitions in the
# numbering part without any paragraph in the body referencing them.
$lastPara = $d.Paragraphs.Last
$insPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insPoint.InsertBefore("T1`rT2`rT3`rT4`rT5`r")

$gallery = $word.ListGalleries.Item(1)
$template = $gallery.ListTemplates.Item(1)

$count = $d.Paragraphs.Count
$firstTemp = $d.Paragraphs.Item($count - 5)
$lastTemp = $d.Paragraphs.Item($count - 1)
for ($i = $count - 5; $i -le $count - 1; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.ListFormat.ApplyListTemplateWithLevel($template)
}
$delRange = $d.Range($firstTemp.Range.Start, $lastTemp.Range.End)
$delRange.Delete()

# --- Step 2: insert the real content before the trailing empty paragraph:
# two blank paragraphs, the "Things to look..." heading line, and the
# twenty bulleted list items.
$lastPara = $d.Paragraphs.Last
$insPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$items = @(
  "Response length vs energy",
  "Prompt length vs efficiency  ",
  "Model efficiency comparison",
  "Time-to-first-token analysis",
  "Question type energy cost",
  "Domain-specific energy use",
  "Code vs text energy difference",
  "Vocabulary complexity impact",
  "Sentence structure efficiency",
  "Conversation depth cost",
  "Readability score correlation",
  "Special content energy cost",
  "Named entity density",
  "Quality vs energy tradeoffs",
  "Time-of-day performance",
  "Tokenization efficiency",
  "Error pattern analysis",
  "Optimal response length",
  "Model-specific optimization",
  "Batch processing efficiency"
)

$text = "`r`rThings to look for correlation with:`r"
foreach ($item in $items) {
    $text += $item + "`r"
}
$insPoint.InsertBefore($text)

# --- Step 3: style the 20 newly-added item paragraphs as a bulleted
# "List Paragraph" list, all sharing one list id.
$count = $d.Paragraphs.Count
$firstItemPara = $d.Paragraphs.Item($count - 20)
$lastItemPara = $d.Paragraphs.Item($count - 1)
$listRange = $d.Range($firstItemPara.Range.Start, $lastItemPara.Range.End)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyListTemplateWithLevel($template)

Write-Output ("Done. Paragraph count: " + $d.Paragraphs.Count)
